# "remove product status from the template upload"
#
# The "Product Properties" sheet has a "Status" column (column E) whose
# header cell reads "Status" and whose data rows carry the fixed value
# "PUBLISHED". That whole column is being removed from the template so the
# uploader no longer has to (or can) set/see a product status.
#
# Column E sits immediately to the left of two merged header bands
# (old F1:M1 = "Dimensions", old N1:AY1 = "MDF raw"). Deleting a column
# that's adjacent to a merge makes Excel drop the stale merge, so we
# unmerge those two ranges first and then delete the entire column,
# exactly like a user would via Home > Delete > Delete Sheet Columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Properties")

$ws.Range("F1:M1").UnMerge()
$ws.Range("N1:AY1").UnMerge()

$ws.Columns.Item(5).Delete()

$ws.Activate() | Out-Null
$ws.Range("E21").Select() | Out-Null
